$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "325.76").
# Excel would auto-convert these to the Number type on assignment, but the
# source data keeps them as plain text (e.g. trailing zeros like "0.4590"
# must be preserved). Force text formatting, assign, then restore the
# original "Normal" cell style so no formatting side effects remain.
$textCells = @(
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D22",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D5").Value = '325.76'
$ws.Range("D7").Value = '0.4590'
$ws.Range("D8").Value = '0.3820'
$ws.Range("D9").Value = '0.07752'
$ws.Range("D10").Value = '0.9796'
$ws.Range("D11").Value = '22.61'
$ws.Range("D13").Value = '5.709'
$ws.Range("D14").Value = '6.981'
$ws.Range("D15").Value = '0.07006'
$ws.Range("D16").Value = '84.87'
$ws.Range("D18").Value = '0.000009507'
$ws.Range("D19").Value = '16.68'
$ws.Range("D22").Value = '5.357'
$ws.Range("D26").Value = '158.06'
$ws.Range("D27").Value = '19.03'
$ws.Range("D28").Value = '5.625'
$ws.Range("D29").Value = '117.65'
$ws.Range("D30").Value = '1.834'
$ws.Range("D31").Value = '0.09327'
$ws.Range("D32").Value = '0.8621'
$ws.Range("D33").Value = '5.104'
$ws.Range("D34").Value = '1.245'
$ws.Range("D35").Value = '3.016'
$ws.Range("D36").Value = '0.05702'
$ws.Range("D37").Value = '1.154'
$ws.Range("D38").Value = '1.004'
$ws.Range("D39").Value = '0.02054'
$ws.Range("D40").Value = '3.108'
$ws.Range("D41").Value = '7.456'
$ws.Range("D42").Value = '0.5512'
$ws.Range("D43").Value = '0.1754'
$ws.Range("D44").Value = '9.343'
$ws.Range("D45").Value = '0.000002804'
$ws.Range("D46").Value = '2.181'
$ws.Range("D47").Value = '0.5185'
$ws.Range("D48").Value = '11.25'
$ws.Range("D49").Value = '0.06938'
$ws.Range("D50").Value = '110.88'
$ws.Range("D51").Value = '1.762'

foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}

# Remaining cells are unambiguously text already (contain letters, "%",
# multiple "." separators, URLs, etc.) so a plain assignment is safe.
$ws.Range("D2").Value = '29.024.82'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.924.68'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("D12").Value = '1.969.40'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '29.051.74'
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = '2.174.25'
$ws.Range("E24").Value = '  -0.97%  '
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("E40").Value = '  +13.74%  '
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("E45").Value = '  +9.45%  '
$ws.Range("E46").Value = '  +4.26%  '
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("E51").Value = '  -0.84%  '
